# Updated cryptos list on Tue Oct 10 04:10:21 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns on Sheet1 with the
# latest scraped values. Cells are plain text (e.g. "27.694.46", which is
# not a valid number, and "  -0.51%  " which carries intentional padding
# spaces), so every write is forced to Text via a leading apostrophe and the
# cell's style is reset to "Normal" afterwards so no numeric coercion or
# stray number-format/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "27.694.46";  "E2"  = "  -0.51%  "
    "D3"  = "1.589.12";   "E3"  = "  -2.45%  "
    "E4"  = "  +0.59%  "
    "D5"  = "207.15";     "E5"  = "  -1.97%  "
    "E7"  = "  +0.59%  "
    "D8"  = "22.23";      "E8"  = "  -4.20%  "
    "E9"  = "  -1.84%  "
    "E10" = "  -2.64%  "
    "E11" = "  -1.43%  "
    "D12" = "1.814.96";   "E12" = "  -2.43%  "
    "D13" = "1.582.68";   "E13" = "  -2.84%  "
    "E14" = "  -3.84%  "
    "E15" = "  -4.24%  "
    "D16" = "27.693.62"
    "E17" = "  -2.10%  "
    "D18" = "219.84";     "E18" = "  -3.78%  "
    "D19" = "0.0₃0696";   "E19" = "  -3.07%  "
    "D20" = "7.32";       "E20" = "  -3.86%  "
    "E21" = "  +0.57%  "
    "E22" = "  -4.73%  "
    "D23" = "9.60";       "E23" = "  -3.36%  "
    "E24" = "  -3.76%  "
    "D25" = "153.70";     "E25" = "  -0.95%  "
    "D26" = "6.86";       "E26" = "  -0.70%  "
    "E27" = "  +0.55%  "
    "D28" = "15.13";      "E28" = "  -2.06%  "
    "E29" = "  -4.50%  "
    "E30" = "  -2.09%  "
    "D31" = "0.0467";     "E31" = "  -2.62%  "
    "E32" = "  -4.96%  "
    "D33" = "1.371.46";   "E33" = "  -3.02%  "
    "E34" = "  -5.49%  "
    "E35" = "  -4.34%  "
    "D36" = "0.981";      "E36" = "  -1.77%  "
    "E37" = "  -0.84%  "
    "E38" = "  -1.30%  "
    "D39" = "0.538";      "E39" = "  -2.67%  "
    "E40" = "  -3.03%  "
    "E41" = "  +0.61%  "
    "D42" = "0.976";      "E42" = "  -2.60%  "
    "D43" = "64.17";      "E43" = "  -2.29%  "
    "E44" = "  +2.48%  "
    "E45" = "  -3.33%  "
    "E46" = "  -4.79%  "
    "D47" = "1.726.01";   "E47" = "  -2.41%  "
    "D48" = "87.84"
    "E49" = "  +10.84%  "
    "E50" = "  -3.88%  "
    "D51" = "0.0495";     "E51" = "  -1.53%  "
}

foreach ($address in $updates.Keys) {
    $range = $ws.Range($address)
    $range.Value = "'" + $updates[$address]
    $range.Style = "Normal"
}
